$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Row, $Col, $Val)
    $c = $ws.Cells.Item($Row, $Col)
    $c.Value = "'" + $Val
    $c.Style = "Normal"
}

Set-TextValue 2 4 '64.592.95'
Set-TextValue 2 5 '  -0.34%  '
Set-TextValue 3 4 '3.421.45'
Set-TextValue 3 5 '  -0.98%  '
Set-TextValue 4 5 '  -0.04%  '
Set-TextValue 5 4 '573.02'
Set-TextValue 5 5 '  -0.47%  '
Set-TextValue 6 4 '156.89'
Set-TextValue 6 5 '  -2.03%  '
Set-TextValue 7 4 '0.627'
Set-TextValue 7 5 '  +7.63%  '
Set-TextValue 8 5 '  +0.02%  '
Set-TextValue 9 4 '3.426.78'
Set-TextValue 9 5 '  -0.84%  '
Set-TextValue 10 4 '7.16'
Set-TextValue 10 5 '  -2.76%  '
Set-TextValue 11 5 '  -1.96%  '
Set-TextValue 12 5 '  +0.72%  '
Set-TextValue 13 4 '4.014.63'
Set-TextValue 13 5 '  -0.96%  '
Set-TextValue 14 5 '  +0.43%  '
Set-TextValue 15 5 '  -3.30%  '
Set-TextValue 16 4 '27.89'
Set-TextValue 16 5 '  -0.79%  '
Set-TextValue 17 4 '64.600.72'
Set-TextValue 17 5 '  -0.43%  '
Set-TextValue 18 4 '3.448.98'
Set-TextValue 18 5 '  -0.36%  '
Set-TextValue 19 4 '6.36'
Set-TextValue 19 5 '  -0.03%  '
Set-TextValue 20 4 '13.98'
Set-TextValue 20 5 '  -2.22%  '
Set-TextValue 21 4 '377.95'
Set-TextValue 21 5 '  -2.53%  '
Set-TextValue 22 4 '8.03'
Set-TextValue 23 5 '  +1.00%  '
Set-TextValue 24 5 '  -0.02%  '
Set-TextValue 25 4 '72.53'
Set-TextValue 25 5 '  -0.93%  '
Set-TextValue 26 5 '  -4.34%  '
Set-TextValue 27 4 '10.34'
Set-TextValue 27 5 '  +6.97%  '
Set-TextValue 28 5 '  -1.71%  '
Set-TextValue 29 5 '  +0.44%  '
Set-TextValue 30 4 '1.49'
Set-TextValue 30 5 '  +4.08%  '
Set-TextValue 31 5 '  -0.61%  '
Set-TextValue 32 4 '2.02'
Set-TextValue 32 5 '  -1.00%  '
Set-TextValue 33 5 '  -2.32%  '
Set-TextValue 34 5 '  +1.69%  '
Set-TextValue 35 5 '  +6.87%  '
Set-TextValue 36 4 '159.55'
Set-TextValue 36 5 '  -2.29%  '
Set-TextValue 37 5 '  -0.19%  '
Set-TextValue 38 4 '6.96'
Set-TextValue 38 5 '  +6.52%  '
Set-TextValue 39 4 '0.0766'
Set-TextValue 39 5 '  +0.16%  '
Set-TextValue 40 4 '26.86'
Set-TextValue 40 5 '  -1.46%  '
Set-TextValue 41 4 '2.884.50'
Set-TextValue 41 5 '  -4.39%  '
Set-TextValue 42 4 '4.63'
Set-TextValue 42 5 '  +1.60%  '
Set-TextValue 43 4 '26.67'
Set-TextValue 43 5 '  +9.15%  '
Set-TextValue 44 5 '  +0.77%  '
Set-TextValue 45 5 '  -0.04%  '
Set-TextValue 46 4 '0.773'
Set-TextValue 46 5 '  -0.24%  '
Set-TextValue 47 4 '322.00'
Set-TextValue 47 5 '  +6.03%  '
Set-TextValue 48 5 '  -0.07%  '
Set-TextValue 49 2 'Stellar'
Set-TextValue 49 3 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 49 4 '0.109'
Set-TextValue 49 5 '  +2.45%  '
Set-TextValue 50 2 'dogwifhat'
Set-TextValue 50 3 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 50 4 '2.19'
Set-TextValue 50 5 '  +0.92%  '
Set-TextValue 51 5 '  -1.31%  '
